# Apply the "etiquetas 302 - HE aguja trozos H 2C - jabat" update:
# Row 9 corresponds to item "2C" / Aguja Trozos "H". A new label number
# (302) was assigned and its translation status flipped from "CERO" to "OK".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "etis numero" (label number) column M, row 9
$ws.Range("M9").Value = 302

# "TRADUCCION" (translation status) column N, row 9
$ws.Range("N9").Value = "OK"

# Row was resized (manually, to better show the now-longer content)
$ws.Rows(9).RowHeight = 93.75

# Move the active selection to N10, as recorded in the saved view state
$ws.Activate()
$ws.Range("N10").Select()
